# control-flow-while-loop.pptx edit:
#   - slide 1 "TextBox 15": rename C#-style Console API calls to the
#     lowercase/underscore style used elsewhere in the deck
#       "using static "         -> "" (both occurrences)
#       "Write(\"Count to: \");" -> "write(\"Count to: \");"
#       "target = ToInt32("     -> "target = to_integer("
#       "ReadLine"              -> "read_line"
#   - slides 2-11 "TextBox 15": "    WriteLine(" -> "    write_line("

$p = $ppt.ActivePresentation

# --- Slide 1 -----------------------------------------------------------
$s1 = $p.Slides.Item(1)
$tb1 = $s1.Shapes.Item(1).GroupItems.Item(1)
$tr1 = $tb1.TextFrame.TextRange

# Apply edits from the highest character offset to the lowest so that
# earlier (lower-offset) text positions are not shifted by the edits
# that follow them.
$tr1.Characters(111, 8).Text = "read_line"
$tr1.Characters(94, 17).Text = "target = to_integer("
$tr1.Characters(73, 20).Text = 'write("Count to: ");'
$tr1.Characters(30, 13).Text = ""
$tr1.Characters(1, 13).Text = ""

# --- Slides 2-11 ---------------------------------------------------------
for ($i = 2; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    $tb = $slide.Shapes.Item(1).GroupItems.Item(1)
    $tr = $tb.TextFrame.TextRange

    $idx = $tr.Text.IndexOf("    WriteLine(")
    if ($idx -ge 0) {
        $tr.Characters($idx + 1, 14).Text = "    write_line("
    }
}
